$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.599.21"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "2.295.59"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "155.79"
$ws.Range("E5").Value = "  +15,460.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.75"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "96.66"
$ws.Range("E7").Value = "  +4.88%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "35.74"
$ws.Range("E11").Value = "  +8.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0809"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "2.653.11"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.57"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "2.315.17"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("E18").Value = "  +5.07%  "
$ws.Range("D19").Value = "42.507.49"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.33"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "244.05"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.97"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.27"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.81"
$ws.Range("E29").Value = "  +7.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.71"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.11"
$ws.Range("E31").Value = "  -8.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.47"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.38"
$ws.Range("E33").Value = "  +4.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0758"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  +5.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.34"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.84"
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  +7.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.84"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").Value = "2.022.77"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("E45").Value = "  +11.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0285"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.01"
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.81"
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.69"
$ws.Range("E51").Value = "  +1.38%  "
